$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Point Masses")
$originalActive = $wb.ActiveSheet

$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 1

$ws.Range("B3").Select()

$originalActive.Activate()
